$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, shifting existing rows 15..125 down to 16..126
$ws.Rows(15).EntireRow.Insert()

# Populate the newly inserted row 15 with the new weekly data point.
# Fixed/template columns shared by every data row in this sheet:
$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(15, 3).Value = "Los Lagos"
$ws.Cells.Item(15, 4).Value = 44490
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 100112039
$ws.Cells.Item(15, 7).Value = "Ciboulette"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 2500
$ws.Cells.Item(15, 12).Value = 2500
$ws.Cells.Item(15, 13).Value = 2500
$ws.Cells.Item(15, 14).Value = "$/docena de atados"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 833
$ws.Cells.Item(15, 17).Value = 3
$ws.Cells.Item(15, 18).Value = "Hortaliza"
